$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.725.55"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.339.93"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Formula = "'513.84"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Formula = "'133.86"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Formula = "'0.533"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Formula = "'0.101"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Formula = "'5.31"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Formula = "'0.340"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Formula = "'23.84"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "2.757.31"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "56.698.29"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Formula = "'0.0000133"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "2.341.14"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Formula = "'10.41"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Formula = "'326.21"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Formula = "'4.18"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").Formula = "'6.67"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Formula = "'61.16"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Formula = "'8.76"
$ws.Range("E24").Value = "  +13.21%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Formula = "'0.165"
$ws.Range("E25").Value = "  +4.16%  "
$ws.Range("D26").Formula = "'1.00"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  +7.92%  "
$ws.Range("D28").Formula = "'168.45"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "0.0₃0727"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Formula = "'1.67"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Formula = "'6.16"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Formula = "'18.41"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Formula = "'0.997"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Formula = "'1.27"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").Formula = "'3.99"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Formula = "'0.887"
$ws.Range("E37").Value = "  -6.14%  "
$ws.Range("D38").Formula = "'1.56"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").Formula = "'38.68"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D40").Formula = "'149.79"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Formula = "'3.60"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").Formula = "'280.87"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").Formula = "'5.12"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Formula = "'0.0926"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Formula = "'0.0499"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Formula = "'0.558"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Formula = "'18.29"
$ws.Range("E48").Value = "  +7.40%  "
$ws.Range("D49").Formula = "'0.0215"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Formula = "'17.08"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").Formula = "'10.99"
$ws.Range("E51").Value = "  +1.19%  "
